# Refresh the NSE ticker lists (columns A-F) on the active sheet.
# Column A keeps its running index (0,1,2,...), column B is the new
# "Buying Opportunity" list, column C the new "support Zone" list (now much
# longer, rows 2-41), and columns D/E/F hold the handful of remaining
# "long buildup" / "Short buildup" / "FII ENTERING" tickers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (index 0)
$ws.Cells.Item(2, 1).Value = 0   # A2
$ws.Cells.Item(2, 2).Value = "NSE:ABFRL"   # B2
$ws.Cells.Item(2, 3).Value = "NSE:ACE"   # C2
$ws.Cells.Item(2, 4).Value = ""   # D2 (blank)
$ws.Cells.Item(2, 5).Value = "NSE:BHARATFORG"   # E2
$ws.Cells.Item(2, 6).Value = "NSE:ABFRL"   # F2

# Row 3 (index 1)
$ws.Cells.Item(3, 1).Value = 1   # A3
$ws.Cells.Item(3, 2).Value = "NSE:AFFLE"   # B3
$ws.Cells.Item(3, 3).Value = "NSE:ADVANIHOTR"   # C3
$ws.Cells.Item(3, 4).Value = ""   # D3 (blank)
$ws.Cells.Item(3, 5).Value = "NSE:DEEPAKNTR"   # E3
$ws.Cells.Item(3, 6).Value = "NSE:AUBANK"   # F3

# Row 4 (index 2)
$ws.Cells.Item(4, 1).Value = 2   # A4
$ws.Cells.Item(4, 2).Value = "NSE:ALMONDZ"   # B4
$ws.Cells.Item(4, 3).Value = "NSE:ALANKIT"   # C4
$ws.Cells.Item(4, 4).Value = ""   # D4 (blank)
$ws.Cells.Item(4, 5).Value = ""   # E4 (blank)
$ws.Cells.Item(4, 6).Value = ""   # F4 (blank)

# Row 5 (index 3)
$ws.Cells.Item(5, 1).Value = 3   # A5
$ws.Cells.Item(5, 2).Value = "NSE:ARENTERP"   # B5
$ws.Cells.Item(5, 3).Value = "NSE:ARMANFIN"   # C5
$ws.Cells.Item(5, 4).Value = ""   # D5 (blank)
$ws.Cells.Item(5, 5).Value = ""   # E5 (blank)
$ws.Cells.Item(5, 6).Value = ""   # F5 (blank)

# Row 6 (index 4)
$ws.Cells.Item(6, 1).Value = 4   # A6
$ws.Cells.Item(6, 2).Value = "NSE:AUBANK"   # B6
$ws.Cells.Item(6, 3).Value = "NSE:ASAHISONG"   # C6
$ws.Cells.Item(6, 4).Value = ""   # D6 (blank)
$ws.Cells.Item(6, 5).Value = ""   # E6 (blank)
$ws.Cells.Item(6, 6).Value = ""   # F6 (blank)

# Row 7 (index 5)
$ws.Cells.Item(7, 1).Value = 5   # A7
$ws.Cells.Item(7, 2).Value = "NSE:FIEMIND"   # B7
$ws.Cells.Item(7, 3).Value = "NSE:AUTOIND"   # C7
$ws.Cells.Item(7, 4).Value = ""   # D7 (blank)
$ws.Cells.Item(7, 5).Value = ""   # E7 (blank)
$ws.Cells.Item(7, 6).Value = ""   # F7 (blank)

# Row 8 (index 6)
$ws.Cells.Item(8, 1).Value = 6   # A8
$ws.Cells.Item(8, 2).Value = "NSE:KEYFINSERV"   # B8
$ws.Cells.Item(8, 3).Value = "NSE:AVANTIFEED"   # C8
$ws.Cells.Item(8, 4).Value = ""   # D8 (blank)
$ws.Cells.Item(8, 5).Value = ""   # E8 (blank)
$ws.Cells.Item(8, 6).Value = ""   # F8 (blank)

# Row 9 (index 7)
$ws.Cells.Item(9, 1).Value = 7   # A9
$ws.Cells.Item(9, 2).Value = "NSE:LUMAXIND"   # B9
$ws.Cells.Item(9, 3).Value = "NSE:AXISCETF"   # C9
$ws.Cells.Item(9, 4).Value = ""   # D9 (blank)
$ws.Cells.Item(9, 5).Value = ""   # E9 (blank)
$ws.Cells.Item(9, 6).Value = ""   # F9 (blank)

# Row 10 (index 8)
$ws.Cells.Item(10, 1).Value = 8   # A10
$ws.Cells.Item(10, 2).Value = "NSE:MOLDTECH"   # B10
$ws.Cells.Item(10, 3).Value = "NSE:BARBEQUE"   # C10
$ws.Cells.Item(10, 4).Value = ""   # D10 (blank)
$ws.Cells.Item(10, 5).Value = ""   # E10 (blank)
$ws.Cells.Item(10, 6).Value = ""   # F10 (blank)

# Row 11 (index 9)
$ws.Cells.Item(11, 1).Value = 9   # A11
$ws.Cells.Item(11, 2).Value = "NSE:S&SPOWER"   # B11
$ws.Cells.Item(11, 3).Value = "NSE:BECTORFOOD"   # C11
$ws.Cells.Item(11, 4).Value = ""   # D11 (blank)
$ws.Cells.Item(11, 5).Value = ""   # E11 (blank)
$ws.Cells.Item(11, 6).Value = ""   # F11 (blank)

# Row 12 (index 10)
$ws.Cells.Item(12, 1).Value = 10   # A12
$ws.Cells.Item(12, 2).Value = ""   # B12 (blank)
$ws.Cells.Item(12, 3).Value = "NSE:CENTURYPLY"   # C12
$ws.Cells.Item(12, 4).Value = ""   # D12 (blank)
$ws.Cells.Item(12, 5).Value = ""   # E12 (blank)
$ws.Cells.Item(12, 6).Value = ""   # F12 (blank)

# Row 13 (index 11)
$ws.Cells.Item(13, 1).Value = 11   # A13
$ws.Cells.Item(13, 2).Value = ""   # B13 (blank)
$ws.Cells.Item(13, 3).Value = "NSE:CHOLAFIN"   # C13
$ws.Cells.Item(13, 4).Value = ""   # D13 (blank)
$ws.Cells.Item(13, 5).Value = ""   # E13 (blank)
$ws.Cells.Item(13, 6).Value = ""   # F13 (blank)

# Row 14 (index 12)
$ws.Cells.Item(14, 1).Value = 12   # A14
$ws.Cells.Item(14, 2).Value = ""   # B14 (blank)
$ws.Cells.Item(14, 3).Value = "NSE:CONSUMBEES"   # C14
$ws.Cells.Item(14, 4).Value = ""   # D14 (blank)
$ws.Cells.Item(14, 5).Value = ""   # E14 (blank)
$ws.Cells.Item(14, 6).Value = ""   # F14 (blank)

# Row 15 (index 13)
$ws.Cells.Item(15, 1).Value = 13   # A15
$ws.Cells.Item(15, 2).Value = ""   # B15 (blank)
$ws.Cells.Item(15, 3).Value = "NSE:ERIS"   # C15
$ws.Cells.Item(15, 4).Value = ""   # D15 (blank)
$ws.Cells.Item(15, 5).Value = ""   # E15 (blank)
$ws.Cells.Item(15, 6).Value = ""   # F15 (blank)

# Row 16 (index 14)
$ws.Cells.Item(16, 1).Value = 14   # A16
$ws.Cells.Item(16, 2).Value = ""   # B16 (blank)
$ws.Cells.Item(16, 3).Value = "NSE:ESTER"   # C16
$ws.Cells.Item(16, 4).Value = ""   # D16 (blank)
$ws.Cells.Item(16, 5).Value = ""   # E16 (blank)
$ws.Cells.Item(16, 6).Value = ""   # F16 (blank)

# Row 17 (index 15)
$ws.Cells.Item(17, 1).Value = 15   # A17
$ws.Cells.Item(17, 2).Value = ""   # B17 (blank)
$ws.Cells.Item(17, 3).Value = "NSE:ETHOSLTD"   # C17
$ws.Cells.Item(17, 4).Value = ""   # D17 (blank)
$ws.Cells.Item(17, 5).Value = ""   # E17 (blank)
$ws.Cells.Item(17, 6).Value = ""   # F17 (blank)

# Row 18 (index 16)
$ws.Cells.Item(18, 1).Value = 16   # A18
$ws.Cells.Item(18, 2).Value = ""   # B18 (blank)
$ws.Cells.Item(18, 3).Value = "NSE:EXXARO"   # C18
$ws.Cells.Item(18, 4).Value = ""   # D18 (blank)
$ws.Cells.Item(18, 5).Value = ""   # E18 (blank)
$ws.Cells.Item(18, 6).Value = ""   # F18 (blank)

# Row 19 (index 17)
$ws.Cells.Item(19, 1).Value = 17   # A19
$ws.Cells.Item(19, 2).Value = ""   # B19 (blank)
$ws.Cells.Item(19, 3).Value = "NSE:FINOPB"   # C19
$ws.Cells.Item(19, 4).Value = ""   # D19 (blank)
$ws.Cells.Item(19, 5).Value = ""   # E19 (blank)
$ws.Cells.Item(19, 6).Value = ""   # F19 (blank)

# Row 20 (index 18)
$ws.Cells.Item(20, 1).Value = 18   # A20
$ws.Cells.Item(20, 2).Value = ""   # B20 (blank)
$ws.Cells.Item(20, 3).Value = "NSE:GALAXYSURF"   # C20
$ws.Cells.Item(20, 4).Value = ""   # D20 (blank)
$ws.Cells.Item(20, 5).Value = ""   # E20 (blank)
$ws.Cells.Item(20, 6).Value = ""   # F20 (blank)

# Row 21 (index 19)
$ws.Cells.Item(21, 1).Value = 19   # A21
$ws.Cells.Item(21, 2).Value = ""   # B21 (blank)
$ws.Cells.Item(21, 3).Value = "NSE:GOKEX"   # C21
$ws.Cells.Item(21, 4).Value = ""   # D21 (blank)
$ws.Cells.Item(21, 5).Value = ""   # E21 (blank)
$ws.Cells.Item(21, 6).Value = ""   # F21 (blank)

# Row 22 (index 20)
$ws.Cells.Item(22, 1).Value = 20   # A22
$ws.Cells.Item(22, 2).Value = ""   # B22 (blank)
$ws.Cells.Item(22, 3).Value = "NSE:HERCULES"   # C22
$ws.Cells.Item(22, 4).Value = ""   # D22 (blank)
$ws.Cells.Item(22, 5).Value = ""   # E22 (blank)
$ws.Cells.Item(22, 6).Value = ""   # F22 (blank)

# Row 23 (index 21)
$ws.Cells.Item(23, 1).Value = 21   # A23
$ws.Cells.Item(23, 2).Value = ""   # B23 (blank)
$ws.Cells.Item(23, 3).Value = "NSE:HPAL"   # C23
$ws.Cells.Item(23, 4).Value = ""   # D23 (blank)
$ws.Cells.Item(23, 5).Value = ""   # E23 (blank)
$ws.Cells.Item(23, 6).Value = ""   # F23 (blank)

# Row 24 (index 22)
$ws.Cells.Item(24, 1).Value = 22   # A24
$ws.Cells.Item(24, 2).Value = ""   # B24 (blank)
$ws.Cells.Item(24, 3).Value = "NSE:INDIAMART"   # C24
$ws.Cells.Item(24, 4).Value = ""   # D24 (blank)
$ws.Cells.Item(24, 5).Value = ""   # E24 (blank)
$ws.Cells.Item(24, 6).Value = ""   # F24 (blank)

# Row 25 (index 23)
$ws.Cells.Item(25, 1).Value = 23   # A25
$ws.Cells.Item(25, 2).Value = ""   # B25 (blank)
$ws.Cells.Item(25, 3).Value = "NSE:JAICORPLTD"   # C25
$ws.Cells.Item(25, 4).Value = ""   # D25 (blank)
$ws.Cells.Item(25, 5).Value = ""   # E25 (blank)
$ws.Cells.Item(25, 6).Value = ""   # F25 (blank)

# Row 26 (index 24)
$ws.Cells.Item(26, 1).Value = 24   # A26
$ws.Cells.Item(26, 2).Value = ""   # B26 (blank)
$ws.Cells.Item(26, 3).Value = "NSE:JASH"   # C26
$ws.Cells.Item(26, 4).Value = ""   # D26 (blank)
$ws.Cells.Item(26, 5).Value = ""   # E26 (blank)
$ws.Cells.Item(26, 6).Value = ""   # F26 (blank)

# Row 27 (index 25)
$ws.Cells.Item(27, 1).Value = 25   # A27
$ws.Cells.Item(27, 2).Value = ""   # B27 (blank)
$ws.Cells.Item(27, 3).Value = "NSE:JYOTISTRUC"   # C27
$ws.Cells.Item(27, 4).Value = ""   # D27 (blank)
$ws.Cells.Item(27, 5).Value = ""   # E27 (blank)
$ws.Cells.Item(27, 6).Value = ""   # F27 (blank)

# Row 28 (index 26)
$ws.Cells.Item(28, 1).Value = 26   # A28
$ws.Cells.Item(28, 2).Value = ""   # B28 (blank)
$ws.Cells.Item(28, 3).Value = "NSE:KANORICHEM"   # C28
$ws.Cells.Item(28, 4).Value = ""   # D28 (blank)
$ws.Cells.Item(28, 5).Value = ""   # E28 (blank)
$ws.Cells.Item(28, 6).Value = ""   # F28 (blank)

# Row 29 (index 27)
$ws.Cells.Item(29, 1).Value = 27   # A29
$ws.Cells.Item(29, 2).Value = ""   # B29 (blank)
$ws.Cells.Item(29, 3).Value = "NSE:KCPSUGIND"   # C29
$ws.Cells.Item(29, 4).Value = ""   # D29 (blank)
$ws.Cells.Item(29, 5).Value = ""   # E29 (blank)
$ws.Cells.Item(29, 6).Value = ""   # F29 (blank)

# Row 30 (index 28)
$ws.Cells.Item(30, 1).Value = 28   # A30
$ws.Cells.Item(30, 2).Value = ""   # B30 (blank)
$ws.Cells.Item(30, 3).Value = "NSE:KEC"   # C30
$ws.Cells.Item(30, 4).Value = ""   # D30 (blank)
$ws.Cells.Item(30, 5).Value = ""   # E30 (blank)
$ws.Cells.Item(30, 6).Value = ""   # F30 (blank)

# Row 31 (index 29)
$ws.Cells.Item(31, 1).Value = 29   # A31
$ws.Cells.Item(31, 2).Value = ""   # B31 (blank)
$ws.Cells.Item(31, 3).Value = "NSE:KEI"   # C31
$ws.Cells.Item(31, 4).Value = ""   # D31 (blank)
$ws.Cells.Item(31, 5).Value = ""   # E31 (blank)
$ws.Cells.Item(31, 6).Value = ""   # F31 (blank)

# Row 32 (index 30)
$ws.Cells.Item(32, 1).Value = 30   # A32
$ws.Cells.Item(32, 2).Value = ""   # B32 (blank)
$ws.Cells.Item(32, 3).Value = "NSE:KHANDSE"   # C32
$ws.Cells.Item(32, 4).Value = ""   # D32 (blank)
$ws.Cells.Item(32, 5).Value = ""   # E32 (blank)
$ws.Cells.Item(32, 6).Value = ""   # F32 (blank)

# Row 33 (index 31)
$ws.Cells.Item(33, 1).Value = 31   # A33
$ws.Cells.Item(33, 2).Value = ""   # B33 (blank)
$ws.Cells.Item(33, 3).Value = "NSE:MENONBE"   # C33
$ws.Cells.Item(33, 4).Value = ""   # D33 (blank)
$ws.Cells.Item(33, 5).Value = ""   # E33 (blank)
$ws.Cells.Item(33, 6).Value = ""   # F33 (blank)

# Row 34 (index 32)
$ws.Cells.Item(34, 1).Value = 32   # A34
$ws.Cells.Item(34, 2).Value = ""   # B34 (blank)
$ws.Cells.Item(34, 3).Value = "NSE:MITTAL"   # C34
$ws.Cells.Item(34, 4).Value = ""   # D34 (blank)
$ws.Cells.Item(34, 5).Value = ""   # E34 (blank)
$ws.Cells.Item(34, 6).Value = ""   # F34 (blank)

# Row 35 (index 33)
$ws.Cells.Item(35, 1).Value = 33   # A35
$ws.Cells.Item(35, 2).Value = ""   # B35 (blank)
$ws.Cells.Item(35, 3).Value = "NSE:MOLDTKPAC"   # C35
$ws.Cells.Item(35, 4).Value = ""   # D35 (blank)
$ws.Cells.Item(35, 5).Value = ""   # E35 (blank)
$ws.Cells.Item(35, 6).Value = ""   # F35 (blank)

# Row 36 (index 34)
$ws.Cells.Item(36, 1).Value = 34   # A36
$ws.Cells.Item(36, 2).Value = ""   # B36 (blank)
$ws.Cells.Item(36, 3).Value = "NSE:MUTHOOTFIN"   # C36
$ws.Cells.Item(36, 4).Value = ""   # D36 (blank)
$ws.Cells.Item(36, 5).Value = ""   # E36 (blank)
$ws.Cells.Item(36, 6).Value = ""   # F36 (blank)

# Row 37 (index 35)
$ws.Cells.Item(37, 1).Value = 35   # A37
$ws.Cells.Item(37, 2).Value = ""   # B37 (blank)
$ws.Cells.Item(37, 3).Value = "NSE:NAGREEKCAP"   # C37
$ws.Cells.Item(37, 4).Value = ""   # D37 (blank)
$ws.Cells.Item(37, 5).Value = ""   # E37 (blank)
$ws.Cells.Item(37, 6).Value = ""   # F37 (blank)

# Row 38 (index 36)
$ws.Cells.Item(38, 1).Value = 36   # A38
$ws.Cells.Item(38, 2).Value = ""   # B38 (blank)
$ws.Cells.Item(38, 3).Value = "NSE:NUVAMA"   # C38
$ws.Cells.Item(38, 4).Value = ""   # D38 (blank)
$ws.Cells.Item(38, 5).Value = ""   # E38 (blank)
$ws.Cells.Item(38, 6).Value = ""   # F38 (blank)

# Row 39 (index 37)
$ws.Cells.Item(39, 1).Value = 37   # A39
$ws.Cells.Item(39, 2).Value = ""   # B39 (blank)
$ws.Cells.Item(39, 3).Value = "NSE:PURVA"   # C39
$ws.Cells.Item(39, 4).Value = ""   # D39 (blank)
$ws.Cells.Item(39, 5).Value = ""   # E39 (blank)
$ws.Cells.Item(39, 6).Value = ""   # F39 (blank)

# Row 40 (index 38)
$ws.Cells.Item(40, 1).Value = 38   # A40
$ws.Cells.Item(40, 2).Value = ""   # B40 (blank)
$ws.Cells.Item(40, 3).Value = "NSE:RUCHIRA"   # C40
$ws.Cells.Item(40, 4).Value = ""   # D40 (blank)
$ws.Cells.Item(40, 5).Value = ""   # E40 (blank)
$ws.Cells.Item(40, 6).Value = ""   # F40 (blank)

# Row 41 (index 39)
$ws.Cells.Item(41, 1).Value = 39   # A41
$ws.Cells.Item(41, 2).Value = ""   # B41 (blank)
$ws.Cells.Item(41, 3).Value = "NSE:SALONA"   # C41
$ws.Cells.Item(41, 4).Value = ""   # D41 (blank)
$ws.Cells.Item(41, 5).Value = ""   # E41 (blank)
$ws.Cells.Item(41, 6).Value = ""   # F41 (blank)

# The sheet previously only had data through row 17; rows 18-41 are brand
# new. Column A there needs the same bold/centered/bordered style used by
# the existing index column, so copy the format from an existing cell and
# paste-special (formats only) onto the new range.
$ws.Range("A2").Copy()
$ws.Range("A18:A41").PasteSpecial(-4122)
$excel.CutCopyMode = 0
